$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "October " / "14" / " - " / "December 4" (4 separate runs) in the
# last Heading1 paragraph get merged into a single run's text.
# A plain Find & Replace across the run boundaries collapses them into one
# run automatically.
# ---------------------------------------------------------------------------
$enDash = [string][char]0x2013
$headingOld = "October 14 " + $enDash + " December 4"
[void]$d.Content.Find.Execute($headingOld, $true, $false, $false, $false, $false, $true, 1, $false, $headingOld, 2)

# ---------------------------------------------------------------------------
# Change 2: split the "Will work on a report for Version 2.0" bullet into
# "Will work on a report for Version " + proofed "2.0", then append two new
# bullet points after it (moving the _GoBack bookmark to the very end).
# ---------------------------------------------------------------------------
$target = $d.Content
[void]$target.Find.Execute("Will work on a report for Version 2.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r = $target.Duplicate

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Will work on a report for Version </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>2.0</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Made a lot of progress with Version 3.0</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Implemented a name screen and displaying it with the player</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@

[void]$r.InsertXML($xml)
